$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the listed rows to match repulled / recalculated data.
$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = -3
$ws.Range("F14").Value = -6
$ws.Range("F17").Value = 5
$ws.Range("F22").Value = 0
$ws.Range("F25").Value = -1
